$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data of row 2 and row 3 (columns B through J); column A (trial index) stays the same
$ws.Range("B2").Value = 8
$ws.Range("C2").Value = 7
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = -5
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 66

$ws.Range("B3").Value = 9
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = 9
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 16

# Match the last-saved selection reflected in the authored workbook
$ws.Range("I15").Select() | Out-Null
